$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.09860586374998093,0.971488893032074,0.02518557384610176,0.9964196085929871),
    @(0.02661536261439323,0.9903790354728699,0.01264103967696428,0.9973146915435791),
    @(0.01646766066551208,0.994846522808075,0.003808316309005022,0.9983888268470764),
    @(0.01311630755662918,0.9951455593109131,0.004480447620153427,0.9983888268470764),
    @(0.01254670321941376,0.9948993325233459,0.005109180696308613,0.9980307817459106),
    @(0.0122664300724864,0.9948289394378662,0.00445084972307086,0.9982098340988159),
    @(0.01278691738843918,0.9945123791694641,0.00468931021168828,0.9982098340988159),
    @(0.01142546813935041,0.9950399994850159,0.004494108259677887,0.9983888268470764),
    @(0.01211865339428186,0.99477618932724,0.004671231843531132,0.9982098340988159),
    @(0.01113103609532118,0.9954973459243774,0.004574078600853682,0.9983888268470764),
    @(0.01219736319035292,0.9946002960205078,0.004766233265399933,0.9982098340988159),
    @(0.01185161899775267,0.9948993325233459,0.004861451219767332,0.9982098340988159),
    @(0.01197303365916014,0.9946706295013428,0.004503251984715462,0.9983888268470764),
    @(0.01186482422053814,0.9946002960205078,0.004596328362822533,0.9983888268470764),
    @(0.0121602788567543,0.9946178793907166,0.004735705442726612,0.9982098340988159),
    @(0.01226719655096531,0.9951103925704956,0.004877014551311731,0.9983888268470764),
    @(0.01121153589338064,0.9948993325233459,0.004190259147435427,0.9983888268470764),
    @(0.01163792703300714,0.9948816895484924,0.003926425706595182,0.9983888268470764),
    @(0.01082005351781845,0.9953742027282715,0.003532412461936474,0.9983888268470764),
    @(0.01082348823547363,0.9953214526176453,0.003465659916400909,0.9983888268470764),
    @(0.01111641898751259,0.9949169158935547,0.003731324104592204,0.9983888268470764),
    @(0.01132191065698862,0.9951807260513306,0.003710911376401782,0.9983888268470764),
    @(0.01050199754536152,0.9949872493743896,0.003639706643298268,0.9983888268470764),
    @(0.01137677021324635,0.9947234392166138,0.00347747397609055,0.9983888268470764),
    @(0.01186009123921394,0.9944244027137756,0.003437234088778496,0.9983888268470764),
    @(0.01105994638055563,0.9949169158935547,0.00333186355419457,0.9983888268470764),
    @(0.01140361092984676,0.994705855846405,0.004107674583792686,0.9983888268470764),
    @(0.01171308010816574,0.99477618932724,0.003441446227952838,0.9983888268470764),
    @(0.01109122298657894,0.9948113560676575,0.003396149957552552,0.9983888268470764),
    @(0.01133560016751289,0.9949169158935547,0.003512562718242407,0.9983888268470764),
    @(0.01069898810237646,0.9953566193580627,0.003510243725031614,0.9983888268470764),
    @(0.01149099785834551,0.994846522808075,0.003312453627586365,0.9983888268470764),
    @(0.01052993629127741,0.9953566193580627,0.003386885160580277,0.9983888268470764),
    @(0.01091468241065741,0.9948289394378662,0.003334768349304795,0.9983888268470764),
    @(0.01085885148495436,0.9950928092002869,0.003491014009341598,0.9983888268470764),
    @(0.01094266586005688,0.9950751662254333,0.003286988940089941,0.9983888268470764),
    @(0.01109958812594414,0.9948993325233459,0.003298672381788492,0.9983888268470764),
    @(0.01118524372577667,0.9949520826339722,0.003441807115450501,0.9983888268470764),
    @(0.01089760567992926,0.9950224161148071,0.003245218424126506,0.9983888268470764),
    @(0.01111915521323681,0.9949344992637634,0.003275980241596699,0.9983888268470764),
    @(0.0106629841029644,0.9952334761619568,0.003240708261728287,0.9983888268470764),
    @(0.01059510000050068,0.9952334761619568,0.003229593392461538,0.9983888268470764),
    @(0.01049361750483513,0.9951807260513306,0.003406606148928404,0.9983888268470764),
    @(0.01156157348304987,0.9949169158935547,0.003254567738622427,0.9983888268470764),
    @(0.01058928575366735,0.9950575828552246,0.003156148362904787,0.9983888268470764),
    @(0.01091843005269766,0.9949696660041809,0.003117673331871629,0.9983888268470764),
    @(0.01047874242067337,0.9952510595321655,0.003077966161072254,0.9983888268470764),
    @(0.01141723897308111,0.9945475459098816,0.00321592390537262,0.9983888268470764),
    @(0.01064102258533239,0.9949344992637634,0.003173034870997071,0.9983888268470764),
    @(0.01126331277191639,0.9948816895484924,0.003117162268608809,0.9983888268470764)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$wb.Save()
